$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.323146939277649
$ws.Range("B1").Value = 1.850236773490906
$ws.Range("C1").Value = 3.315932035446167
$ws.Range("D1").Value = 3.78935980796814
$ws.Range("E1").Value = 1.115804076194763
